# Add a new "2021" column (column O) to the sheet, mirroring the
# formatting already used for column N in each row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4 (header year) ---
$ws.Range("N4").Copy()
$ws.Range("O4").PasteSpecial(-4122)
$ws.Range("O4").Value = 2021

# --- Row 5 (Total) ---
$ws.Range("N5").Copy()
$ws.Range("O5").PasteSpecial(-4122)
$ws.Range("O5").Value = 689
$ws.Range("O5").NumberFormat = "0.0"

# --- Row 6 ---
$ws.Range("N10").Copy()
$ws.Range("O6").PasteSpecial(-4122)
$ws.Range("O6").Value = 94.1

# --- Row 7 ---
$ws.Range("N10").Copy()
$ws.Range("O7").PasteSpecial(-4122)
$ws.Range("O7").Value = 147.1

# --- Row 8 ---
$ws.Range("N10").Copy()
$ws.Range("O8").PasteSpecial(-4122)
$ws.Range("O8").Value = 10.1

# --- Row 9 (dash / no data) ---
$ws.Range("N9").Copy()
$ws.Range("O9").PasteSpecial(-4122)
$ws.Range("O9").Value = "-"
$ws.Range("O9").NumberFormat = "0.0"

# --- Row 10 ---
$ws.Range("N10").Copy()
$ws.Range("O10").PasteSpecial(-4122)
$ws.Range("O10").Value = 82.1

# --- Row 11 ---
$ws.Range("N10").Copy()
$ws.Range("O11").PasteSpecial(-4122)
$ws.Range("O11").Value = 145.3

# --- Row 12 ---
$ws.Range("N10").Copy()
$ws.Range("O12").PasteSpecial(-4122)
$ws.Range("O12").Value = 98.8

# --- Row 13 ---
$ws.Range("N10").Copy()
$ws.Range("O13").PasteSpecial(-4122)
$ws.Range("O13").Value = 98.7

# --- Row 14 ---
$ws.Range("N10").Copy()
$ws.Range("O14").PasteSpecial(-4122)
$ws.Range("O14").Value = 1.8

# --- Row 15 (dash / no data) ---
$ws.Range("N9").Copy()
$ws.Range("O15").PasteSpecial(-4122)
$ws.Range("O15").Value = "-"
$ws.Range("O15").NumberFormat = "0.0"

# --- Row 16 ---
$ws.Range("N16").Copy()
$ws.Range("O16").PasteSpecial(-4122)
$ws.Range("O16").Value = 10.9
$ws.Range("O16").NumberFormat = "0.0"

# Clear clipboard marching-ants selection state
$excel.CutCopyMode = $false

# Update the active selection the way the saved workbook recorded it
$ws.Range("P5").Select()

Write-Output "done"
